# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
# to match the newly scraped numbers from the upstream data source.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1557
$ws1.Range("F3").Value  = 8878
$ws1.Range("F5").Value  = 499
$ws1.Range("F7").Value  = 325
$ws1.Range("F9").Value  = 36
$ws1.Range("F10").Value = 49
$ws1.Range("F11").Value = 3752
$ws1.Range("F15").Value = 4037
$ws1.Range("F17").Value = 1132
$ws1.Range("F20").Value = 2543
$ws1.Range("F21").Value = 92

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1557
$ws4.Range("F3").Value  = 8878
$ws4.Range("F5").Value  = 499
$ws4.Range("F7").Value  = 325
$ws4.Range("F9").Value  = 36
$ws4.Range("F10").Value = 49
$ws4.Range("F11").Value = 3752
$ws4.Range("F15").Value = 4037
$ws4.Range("F17").Value = 1132
$ws4.Range("F20").Value = 2543
$ws4.Range("F22").Value = 92

$wb.Save()
